# Update column G ("K" = strikeouts) values for rows 2-27 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 0
    3  = 3
    4  = 3
    5  = 2
    6  = 3
    7  = 2
    8  = 4
    9  = 4
    10 = 3
    11 = 4
    12 = 1
    13 = 3
    14 = 3
    15 = 6
    16 = 8
    17 = 8
    18 = 4
    19 = 5
    20 = 6
    21 = 3
    22 = 2
    23 = 5
    24 = 5
    25 = 6
    26 = 5
    27 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
